$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2:I3").Value = 0.8220858895705522
$ws.Range("J2:J3").Value = 0.8220858895705522
$ws.Range("K2:K3").Value = 47.5
$ws.Range("L2:L3").Value = 2.914110429447853
$ws.Range("M2:M3").Value = 16.8
$ws.Range("N2:N3").Value = 0.1141304347826087
$ws.Range("O2:O3").Value = 0.3536842105263158
$ws.Range("P2:P3").Value = 16.8
$ws.Range("Q2:Q3").Value = 0.1141304347826087
$ws.Range("R2:R3").Value = 0.3536842105263158
$ws.Range("U2:U3").Value = 7.16
$ws.Range("V2:V3").Value = 0.04864130434782609
$ws.Range("W2:W3").Value = 0.1602564102564103
$ws.Range("X2:X3").Value = 0.08177085369517975
$ws.Range("Y2:Y3").Value = 0.07848555656123052
$ws.Range("Z2:Z3").Value = 0.05691340782122906
$ws.Range("AA2:AA3").Value = 0.04678770949720671
$ws.Range("AB2:AB3").Value = 0.08177085369517975
$ws.Range("AC2:AC3").Value = -0.03498314419797304
$ws.Range("AG2:AG3").Value = -7.16
$ws.Range("AJ2:AJ3").Value = -0.0511282490716938
$ws.Range("AK2:AK3").Value = -0.02364284770836085
$ws.Range("AM2:AM3").Value = -0.9409999999999999
$ws.Range("AQ2:AQ3").Value = -14.24017003188098
